$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential notice date (2021-04-06 -> 2021-04-08) in A40
$ws.Range("A40").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-37
$ws.Range("D2").Value = 0.03111736817043126
$ws.Range("E2").Value = 0.008459271932325763
$ws.Range("D3").Value = 0.03502394212483031
$ws.Range("E3").Value = 0.01303175857805883
$ws.Range("D4").Value = 0.03306408917886941
$ws.Range("E4").Value = 0.01399026763990241
$ws.Range("D5").Value = 0.03879706720061992
$ws.Range("E5").Value = 0.006071251055836724
$ws.Range("D6").Value = 0.01620078850778008
$ws.Range("E6").Value = 0.004454505622900395
$ws.Range("D7").Value = 0.0164953698089658
$ws.Range("E7").Value = 0.006526572473642522
$ws.Range("D8").Value = 0.03228942317077139
$ws.Range("E8").Value = -0.004997581815250784
$ws.Range("D9").Value = 0.03204050788655662
$ws.Range("E9").Value = 0.02020470556958642
$ws.Range("D10").Value = 0.03235602930433265
$ws.Range("E10").Value = 0.01387593923106456
$ws.Range("D11").Value = 0.02933532875687282
$ws.Range("E11").Value = 0.0125099813681131
$ws.Range("D12").Value = 0.01672488295125508
$ws.Range("E12").Value = 0.05581099243120868
$ws.Range("D13").Value = 0.01656280408273121
$ws.Range("E13").Value = 0.01124999999999998
$ws.Range("D14").Value = 0.008182025216869216
$ws.Range("E14").Value = 0.00728744939271242
$ws.Range("D15").Value = 0.008042897662574273
$ws.Range("E15").Value = -0.005648387855966108
$ws.Range("D16").Value = 0.03134487011508191
$ws.Range("E16").Value = 0.01657683771909935
$ws.Range("D17").Value = 0.03207174060282691
$ws.Range("E17").Value = 0.0006307821698907645
$ws.Range("D18").Value = 0.03152540467958368
$ws.Range("E18").Value = 0.01853839398965751
$ws.Range("D19").Value = 0.03333631069454344
$ws.Range("E19").Value = -0.0002235778849531966
$ws.Range("D20").Value = 0.02661500649202767
$ws.Range("E20").Value = 0.007005440773798943
$ws.Range("D21").Value = 0.03090619241837643
$ws.Range("E21").Value = 0.01405604042259978
$ws.Range("D22").Value = 0.03327029608969942
$ws.Range("E22").Value = 0.006272602169808073
$ws.Range("D23").Value = 0.03204855267711109
$ws.Range("E23").Value = -0.01004075364715606
$ws.Range("D24").Value = 0.01723288781362113
$ws.Range("E24").Value = -0.0137302284710018
$ws.Range("D25").Value = 0.01569680604069126
$ws.Range("E25").Value = 0.006104914078987056
$ws.Range("D26").Value = 0.03192598792689889
$ws.Range("E26").Value = 0.001408137552805
$ws.Range("D27").Value = 0.03200927517028633
$ws.Range("E27").Value = -0.01961088688812995
$ws.Range("D28").Value = 0.03252106581644273
$ws.Range("E28").Value = 0.01340536214485799
$ws.Range("D29").Value = 0.0323069324208017
$ws.Range("E29").Value = -0.02094624285923519
$ws.Range("D30").Value = 0.0334651456491584
$ws.Range("E30").Value = 0.01226711917135082
$ws.Range("D31").Value = 0.03223583066898941
$ws.Range("E31").Value = -0.006066522557701792
$ws.Range("D32").Value = 0.03404082139677676
$ws.Range("E32").Value = 0.02952706647760439
$ws.Range("D33").Value = 0.03028011842404919
$ws.Range("E33").Value = -0.002051197899573309
$ws.Range("D34").Value = 0.0459120929173309
$ws.Range("E34").Value = 0.0003478664192950731
$ws.Range("D35").Value = 0.03115948501509877
$ws.Range("E35").Value = 0.00122636029174461
$ws.Range("D36").Value = 0.03386265294714395
$ws.Range("E36").Value = -0.004842259721203113
$ws.Range("E37").Value = 0.005741969258252544
